# "break neighbourhood symettries, about to set small iterations to investigate"
#
# Two data points in the "04 Sep" sheet move (their dependent J-column
# difference formulas - G minus E - recompute automatically), the active
# selection moves off the J column and onto E11, and the three existing
# cellIs conditional-format rules on J6:J22 get their fill/font colours
# reshuffled (the stale duplicate green rule is effectively retired in
# favour of reusing red/yellow).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("04 Sep")

# --- data edits -----------------------------------------------------
# E10: 193 -> 192  (J10 = G10-E10 ripples -11 -> -10)
$ws.Range("E10").Value2 = 192

# E22: 318 -> 302  (J22 = G22-E22 ripples -52 -> -36)
$ws.Range("E22").Value2 = 302

# --- conditional-format colour reshuffle -----------------------------
# Before: lessThan=red(9C0006/FFC7CE), greaterThan=green(006100/C6EFCE),
#         equal=yellow(9C5700/FFEB9C)
# After:  lessThan=red(unchanged), greaterThan=yellow(9C5700/FFEB9C),
#         equal=red(9C0006/FFC7CE)
$fcs = $ws.Range("J6:J22").FormatConditions
for ($i = 1; $i -le $fcs.Count; $i++) {
    $fc = $fcs.Item($i)
    if ($fc.Operator -eq 3) {
        # xlEqual -> red
        $fc.Font.Color = 393372
        $fc.Interior.Color = 13551615
    } elseif ($fc.Operator -eq 5) {
        # xlGreater -> yellow
        $fc.Font.Color = 22428
        $fc.Interior.Color = 10284031
    } elseif ($fc.Operator -eq 6) {
        # xlLess -> red (stays the same colour it already had)
        $fc.Font.Color = 393372
        $fc.Interior.Color = 13551615
    }
}

# --- selection moves to E11 ------------------------------------------
$ws.Activate()
$ws.Range("E11").Select()
